$wb = $excel.ActiveWorkbook

# --- Text / status updates -------------------------------------------------
# "Handed back: in sync with en-US" -> "Ready for handoff"
# "2016-09-06 03:06:13" -> "2016-09-06 03:07:02"
# "2016-09-06 03:06:06" -> "2016-09-06 03:06:57"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 03:07:02"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 03:06:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 03:07:02"

# --- Column width updates ---------------------------------------------------
# Columns previously sized to fit the long status text are narrowed now that
# the status text is shorter ("Ready for handoff").

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
